# Adds three new custom character styles (GaNStyle, GaNParagraph, GaNLinks)
# to the document and applies them to the runs in the relevant paragraphs,
# as described by the commit "Add styles to the new paragraphs".

$d = $word.ActiveDocument

# --- Define the three new character styles -------------------------------

$styleGaNStyle = $d.Styles.Add("GaNStyle", 2)              # wdStyleTypeCharacter
$styleGaNStyle.Font.Name = "Calibri"
$styleGaNStyle.Font.Size = 14                               # w:sz 28 (half-points)

$styleGaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$styleGaNParagraph.Font.Name = "Calibri"
$styleGaNParagraph.Font.Size = 10                            # w:sz 20

$styleGaNLinks = $d.Styles.Add("GaNLinks", 2)
$styleGaNLinks.Font.Name = "Calibri"
$styleGaNLinks.Font.Size = 9.5                                # w:sz 19
$styleGaNLinks.Font.Bold = $true
$styleGaNLinks.Font.Color = 8388608                           # navy -> w:color 000080
$styleGaNLinks.Font.Underline = 1                             # wdUnderlineSingle

# --- Apply GaNStyle to every "2022 Campaign Dates ..." run ---------------

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Text = " 2022 Campaign Dates that use Pegasus constellation: October 8-17, November 7-16,"
$rng.Find.Forward = $true
$rng.Find.Wrap = 1
$n1 = 0
while ($rng.Find.Execute()) {
    $n1 = $n1 + 1
    $rng.Style = "GaNStyle"
    $rng.Collapse(0)
}

# --- Apply GaNParagraph to the "You are participating ..." run -----------

$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$rng2.Find.Text = "You are participating in a global campaign to observe and record the faintest stars visible as a means of measuring light pollution in a given location. By locating and observing the constellation Pegasus constellation in the night sky and comparing it to stellar charts, people from around the world will learn how the lights in their community contribute to light pollution. Your contributions to the online database will document the visible nighttime sky."
$rng2.Find.Forward = $true
$rng2.Find.Wrap = 1
$n2 = 0
while ($rng2.Find.Execute()) {
    $n2 = $n2 + 1
    $rng2.Style = "GaNParagraph"
    $rng2.Collapse(0)
}

# --- Apply GaNLinks to the URL run ----------------------------------------

$rng3 = $d.Content
$rng3.Find.ClearFormatting()
$rng3.Find.Text = "(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
$rng3.Find.Forward = $true
$rng3.Find.Wrap = 1
$n3 = 0
while ($rng3.Find.Execute()) {
    $n3 = $n3 + 1
    $rng3.Style = "GaNLinks"
    $rng3.Collapse(0)
}

Write-Output "GaNStyle applied: $n1 (expected 4)"
Write-Output "GaNParagraph applied: $n2 (expected 1)"
Write-Output "GaNLinks applied: $n3 (expected 1)"
